# Add 2022-Q4 data
# 1) Update the "总计" (summary) sheet: insert a new data row for 2022-Q4 at the
#    top of the data (row 2) and push the existing quarters down by one row,
#    re-numbering the index column (A) sequentially and adding the trailing
#    2021-Q1 row that falls out at the bottom.
$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item("总计")

$quarters = @(
  @("2022-Q4", 5, 0.72),
  @("2022-Q3", 8, 3.25),
  @("2022-Q2", 1, 0.07000000000000001),
  @("2022-Q1", 6, 0.99),
  @("2021-Q4", 2, 1.4),
  @("2021-Q3", 5, 0.11),
  @("2021-Q2", 2, 0.05),
  @("2021-Q1", 3, 0.09)
)

for ($i = 0; $i -lt $quarters.Count; $i++) {
  $row = $i + 2
  $summary.Cells.Item($row, 1).Value = $i
  $summary.Cells.Item($row, 2).Value = $quarters[$i][0]
  $summary.Cells.Item($row, 3).Value = $quarters[$i][1]
  $summary.Cells.Item($row, 4).Value = $quarters[$i][2]
}

# Row 9 is brand new (the table used to end at row 8); copy the index-column
# formatting from the row above so it keeps the same "bold/centered" style
# used by the rest of column A instead of picking up the default style.
$summary.Cells.Item(8, 1).Copy()
$summary.Cells.Item(9, 1).PasteSpecial(-4122)
$summary.Cells.Item(9, 1).Value = $quarters.Count - 1

# 2) Insert a brand-new worksheet named "2022-Q4" right before the existing
#    "2022-Q3" sheet (so the tab order becomes 总计, 2022-Q4, 2022-Q3, ...).
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Count; $c++) {
  $q4.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$funds = @(
  @(0, "002350", "华安安华灵活配置混合A", "28.67", "94.34", "1.72", "0.4931", 8),
  @(1, "009970", "财通内需增长12个月定期开放混合", "8.87", "65.22", "2.30", "0.2040", 9),
  @(2, "016183", "华安安华灵活配置混合C", "0.98", "94.34", "1.72", "0.0169", 8),
  @(3, "519034", "海富通中证500指数增强A", "0.24", "92.50", "1.63", "0.0039", 6),
  @(4, "009004", "海富通中证500指数增强C", "0.04", "92.50", "1.63", "0.0007", 6)
)

foreach ($fund in $funds) {
  $row = [int]$fund[0] + 2
  $q4.Cells.Item($row, 1).Value = $fund[0]
  $q4.Cells.Item($row, 2).Value = $fund[1]
  $q4.Cells.Item($row, 3).Value = $fund[2]
  $q4.Cells.Item($row, 4).Value = $fund[3]
  $q4.Cells.Item($row, 5).Value = $fund[4]
  $q4.Cells.Item($row, 6).Value = $fund[5]
  $q4.Cells.Item($row, 7).Value = $fund[6]
  $q4.Cells.Item($row, 8).Value = $fund[7]
}
